$wb = $excel.ActiveWorkbook

# --- gencode sheet: NCBI/UniProt version bump 2022_03 -> 2022_04 ---
$wsGencode = $wb.Worksheets.Item("gencode")
$wsGencode.Range("E4").Value = "2022_04"

# --- predisposition sheet: 3 curated genes / new date stamp + formatting ---
$wsPredisp = $wb.Worksheets.Item("predisposition")

# New value for E7 (curation/version date stamp)
$wsPredisp.Range("E7").Value = 20221128

# Left-align column E header + data cells (mirrors gencode sheet's column E formatting)
$wsPredisp.Range("E1").HorizontalAlignment = -4131  # xlLeft
$wsPredisp.Range("E2").HorizontalAlignment = -4131  # xlLeft
$wsPredisp.Range("E3").HorizontalAlignment = -4131  # xlLeft
$wsPredisp.Range("E7").HorizontalAlignment = -4131  # xlLeft
$wsPredisp.Range("E8").HorizontalAlignment = -4131  # xlLeft

# --- Switch active tab from gencode to predisposition, updating selections ---
$wsGencode.Range("E4").Select()
$wsPredisp.Activate()
$wsPredisp.Range("C16").Select()
